# AddAppointmentSequenceDiagram.pptx edit:
#   - Bump the auto date placeholder text on the slide master + every slide
#     layout from "11 Nov 2018" to "12 Nov 2018".
#   - Rename the "AddressBook" concept to "HealthBook" in the two labels on
#     the sequence diagram slide (":Address" -> ":Health" and
#     "VersionedAddressBook" -> "VersionedHealthBook").

$p = $ppt.ActivePresentation

# --- 1. Date placeholder: slide master -------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "12 Nov 2018"
    }
}

# --- 2. Date placeholder: every slide layout --------------------------------
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "12 Nov 2018"
        }
    }
}

# --- 3. ":Address" / "BookParser" box -> ":Health" / "BookParser" ----------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text.StartsWith(":Address") -and $tr.Text.Contains("BookParser")) {
            $tr.Characters(1, 8).Text = ":Health"
        }
    }
}

# --- 4. ":VersionedAddressBook" -> ":VersionedHealthBook" (inside group) ---
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Type -eq 6) {
        # msoGroup
        for ($j = 1; $j -le $shp.GroupItems.Count; $j++) {
            $sub = $shp.GroupItems.Item($j)
            if ($sub.HasTextFrame -and $sub.TextFrame.HasText) {
                $subTr = $sub.TextFrame.TextRange
                if ($subTr.Text -eq ":VersionedAddressBook") {
                    $subTr.Characters(2, 20).Text = "VersionedHealthBook"
                }
            }
        }
    }
}
